# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
#
# Change: cell B11 on the active sheet ("Rules") held the text "R40".
# It is retyped as the text "1" (a new shared string), keeping the
# cell's existing number format / style ("s=23", still a text cell).
#
# A plain  $ws.Range("B11").Value = "1"  would get auto-coerced to a
# Number by Excel's type inference (and would also restyle the cell).
# To force the entry to stay Text - exactly like a user typing into a
# Text-formatted cell - we stage the value on a scratch cell that is
# explicitly formatted as Text, copy it, and paste only the value into
# B11 so the destination's existing formatting/style is left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "1"

$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues

$excel.CutCopyMode = 0
$scratch.Clear()
